# Update cryptocurrency price/volume data in cryptos.xlsx (Sheet1)
# Mirrors a scheduled GitHub Actions refresh of coinranking.com data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '70.833.42'
$ws.Range('E2').Value = '  +2.31%  '
# Row 3
$ws.Range('D3').Value = '3.569.66'
$ws.Range('E3').Value = '  +1.52%  '
# Row 4
$ws.Range('E4').Value = '  +0.04%  '
# Row 5
$ws.Range('D5').Value = "'581.90"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.42%  '
# Row 6
$ws.Range('D6').Value = "'188.26"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.92%  '
# Row 7
$ws.Range('E7').Value = '  +2.27%  '
# Row 8
$ws.Range('D8').Value = '3.559.44'
$ws.Range('E8').Value = '  +1.41%  '
# Row 9
$ws.Range('E9').Value = '  -0.07%  '
# Row 10
$ws.Range('D10').Value = "'0.218"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +17.57%  '
# Row 11
$ws.Range('D11').Value = "'0.653"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.28%  '
# Row 12
$ws.Range('D12').Value = "'54.64"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.44%  '
# Row 13
$ws.Range('D13').Value = "'0.0000318"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.44%  '
# Row 14
$ws.Range('D14').Value = "'9.57"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.00%  '
# Row 15
$ws.Range('D15').Value = '4.136.81'
$ws.Range('E15').Value = '  +1.39%  '
# Row 16
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '70.888.96'
$ws.Range('E16').Value = '  +2.54%  '
# Row 17
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = "'19.26"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.29%  '
# Row 18
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').Value = "'12.83"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.93%  '
# Row 19
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.547.26'
$ws.Range('E19').Value = '  +1.07%  '
# Row 20
$ws.Range('D20').Value = "'577.14"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +5.99%  '
# Row 21
$ws.Range('E21').Value = '  +0.69%  '
# Row 22
$ws.Range('E22').Value = '  -1.06%  '
# Row 23
$ws.Range('E23').Value = '  -3.92%  '
# Row 24
$ws.Range('E24').Value = '  +3.27%  '
# Row 25
$ws.Range('D25').Value = "'4.90"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -2.08%  '
# Row 26
$ws.Range('D26').Value = "'94.21"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.26%  '
# Row 27
$ws.Range('E27').Value = '  +1.66%  '
# Row 28
$ws.Range('E28').Value = '  +1.67%  '
# Row 29
$ws.Range('E29').Value = '  +2.26%  '
# Row 30
$ws.Range('D30').Value = "'32.79"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.94%  '
# Row 31
$ws.Range('D31').Value = "'7.25"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.17%  '
# Row 32
$ws.Range('D32').Value = "'12.35"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.30%  '
# Row 33
$ws.Range('E33').Value = '  +2.10%  '
# Row 34
$ws.Range('D34').Value = "'3.81"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +22.95%  '
# Row 35
$ws.Range('D35').Value = "'63.13"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.70%  '
# Row 36
$ws.Range('D36').Value = "'3.31"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +6.54%  '
# Row 37
$ws.Range('D37').Value = "'539.79"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -3.39%  '
# Row 38
$ws.Range('D38').Value = "'0.412"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.55%  '
# Row 39
$ws.Range('D39').Value = '0.0₃0824'
$ws.Range('E39').Value = '  +7.43%  '
# Row 40
$ws.Range('D40').Value = "'38.23"
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.25%  '
# Row 41
$ws.Range('E41').Value = '  +0.02%  '
# Row 42
$ws.Range('D42').Value = '3.626.64'
$ws.Range('E42').Value = '  +10.71%  '
# Row 43
$ws.Range('E43').Value = '  +5.17%  '
# Row 44
$ws.Range('E44').Value = '  +2.46%  '
# Row 45
$ws.Range('E45').Value = '  +5.67%  '
# Row 46
$ws.Range('B46').Value = 'ThetaToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D46').Value = "'2.95"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.62%  '
# Row 47
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = "'3.47"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.39%  '
# Row 48
$ws.Range('D48').Value = "'9.34"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.43%  '
# Row 49
$ws.Range('E49').Value = '  +2.50%  '
# Row 50
$ws.Range('D50').Value = "'0.999"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.15%  '
# Row 51
$ws.Range('E51').Value = '  +5.70%  '

Write-Output "Applied 95 cell updates across rows 2-51"
